# Fruta / hortaliza, semanal
# Insert a new weekly record at row 100 of the Apio / Macroferia Regional de
# Talca sheet. This shifts the existing rows 100-218 down to 101-219 (Excel
# takes care of the shift + the sheet's used-range "dimension") and we only
# need to populate the brand-new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 100..218 down one position, creating a blank row 100.
$ws.Rows(100).Insert()

# Fill in the new row 100 with the new weekly observation.
$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 44789
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = 100112017
$ws.Range("G100").Value = "Apio"
$ws.Range("H100").Value = "Americana (o)"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = 10000
$ws.Range("N100").Value = "$/docena de matas"
$ws.Range("O100").Value = "Provincia del Elquí"
$ws.Range("P100").Value = 1667
$ws.Range("Q100").Value = 6
$ws.Range("R100").Value = "Hortaliza"
